$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update total_risk (R) and total_risk_resp (S) columns with newest airtoxics NATA data
$ws.Range("R2").Value = 44.4444444444444

$ws.Range("S3").Value = 0.5

$ws.Range("R4").Value = 30
$ws.Range("S4").Value = 0.4

$ws.Range("R5").Value = 29.7674418604651
$ws.Range("S5").Value = 0.354651162790698

$ws.Range("R6").Value = 20
$ws.Range("S6").Value = 0.27

$ws.Range("R7").Value = 21.2307692307692
$ws.Range("S7").Value = 0.3

$ws.Range("R9").Value = 65
$ws.Range("S9").Value = 0.4

$ws.Range("R10").Value = 62.5
$ws.Range("S10").Value = 0.4375

$ws.Range("R11").Value = 35
$ws.Range("S11").Value = 0.425

$ws.Range("S12").Value = 0.32

$ws.Range("R13").Value = 47.8947368421053
$ws.Range("S13").Value = 0.3

$ws.Range("R14").Value = 19.0909090909091
$ws.Range("S14").Value = 0.2

$ws.Range("S15").Value = 0.3
